$d = $word.ActiveDocument

# --- Edit 1 -----------------------------------------------------------
# "bouche de leur <tl>fonte</tl> deulx ou" -> "bouche de leur fonte deulx ou"
# (drop the <tl>...</tl> tag-marker runs, merging the plain text around them)
$null = $d.Content.Find.Execute(
    "bouche de leur <tl>fonte</tl> deulx ou", $true, $false, $false, $false, $false,
    $true, 1, $false, "bouche de leur fonte deulx ou", 2)

# --- Edit 2 -------------------------------------------------------------
# "veulent faire courre la <m>fonte</m>." -> "veulent faire courre la fonte."
# Done surgically (delete just the tag-marker runs) so the surrounding
# plain-text runs keep their own formatting/boundaries as much as possible.
$rng = $d.Content
$found = $rng.Find.Execute(
    "veulent faire courre la <m>fonte</m>.", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0)
if ($found) {
    $base = $rng.Start
    # "</m>" sits at relative offset 32..36 within the matched text; remove it
    # first so the trailing "." folds back onto "onte".
    $closeTag = $d.Range($base + 32, $base + 36)
    $closeTag.Delete()
    # "<m>" sits at relative offset 24..27 within the matched text; remove it
    # last.
    $openTag = $d.Range($base + 24, $base + 27)
    $openTag.Delete()
}

# --- Edit 3 -----------------------------------------------------------
# "e aulx aultres <tl>fontes</tl>, affin de mectre" ->
# "e aulx aultres fontes, affin de mectre"
$null = $d.Content.Find.Execute(
    "e aulx aultres <tl>fontes</tl>, affin de mectre", $true, $false, $false, $false, $false,
    $true, 1, $false, "e aulx aultres fontes, affin de mectre", 2)
